$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.971817374229431
$ws.Range("B1").Value = 1.999925494194031
$ws.Range("C1").Value = 8.136432647705078
$ws.Range("D1").Value = 0.9491493701934814
$ws.Range("E1").Value = 0.4557609558105469
